$d = $word.ActiveDocument

# Locate the phrase "Check out the magister branch:" and find the position
# of the erroneous "gi" in "magister" so we can remove it, fixing the typo
# to read "master".
$findRange = $d.Content
$findRange.Find.Execute("Check out the magister branch", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

$giStart = $findRange.Start + "Check out the ma".Length
$giEnd = $giStart + "gi".Length

$giRange = $d.Range($giStart, $giEnd)
$giRange.Text = ""

# Word automatically tracks the location of the most recent edit with the
# hidden "_GoBack" bookmark; move it here (this also removes it from its
# previous location, since _GoBack is a singleton bookmark).
$lastEditRange = $d.Range($giStart, $giStart)
$d.Bookmarks.Add("_GoBack", $lastEditRange)
